# Update "想去人数" (and one "最低票价") figures across the four sheets of the
# workbook to match the regenerated data output.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 567
$ws1.Range("F5").Value  = 6332
$ws1.Range("F6").Value  = 712
$ws1.Range("F7").Value  = 1081
$ws1.Range("F8").Value  = 64
$ws1.Range("F9").Value  = 622
$ws1.Range("F10").Value = 309
$ws1.Range("F12").Value = 673
$ws1.Range("F14").Value = 1141
$ws1.Range("F16").Value = 397
$ws1.Range("F17").Value = 45
$ws1.Range("F18").Value = 15
$ws1.Range("F19").Value = 1410
$ws1.Range("F20").Value = 654
$ws1.Range("F21").Value = 370
$ws1.Range("F22").Value = 389
$ws1.Range("F24").Value = 1062
$ws1.Range("F25").Value = 118
$ws1.Range("F26").Value = 2179
$ws1.Range("F27").Value = 239
$ws1.Range("F28").Value = 84
$ws1.Range("F29").Value = 388
$ws1.Range("F31").Value = 3518

# --- Sheet: 演出 (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value  = 163
$ws2.Range("F7").Value  = 26
$ws2.Range("F19").Value = 312
$ws2.Range("F24").Value = 182
$ws2.Range("F28").Value = 206
$ws2.Range("F29").Value = 29
$ws2.Range("G29").Value = 328
$ws2.Range("F32").Value = 1588

# --- Sheet: 本地生活 (Local Life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F7").Value  = 1566
$ws3.Range("F11").Value = 740

# --- Sheet: 全部类型 (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 1566
$ws4.Range("F8").Value  = 740
$ws4.Range("F9").Value  = 567
$ws4.Range("F12").Value = 6332
$ws4.Range("F14").Value = 26
$ws4.Range("F15").Value = 712
$ws4.Range("F16").Value = 1081
$ws4.Range("F18").Value = 622
$ws4.Range("F20").Value = 673
$ws4.Range("F25").Value = 1141
$ws4.Range("F26").Value = 397
$ws4.Range("F29").Value = 45
$ws4.Range("F30").Value = 15
$ws4.Range("F31").Value = 1410
$ws4.Range("F34").Value = 654
$ws4.Range("F35").Value = 370
$ws4.Range("F36").Value = 389
$ws4.Range("F39").Value = 182
$ws4.Range("F43").Value = 206
$ws4.Range("F45").Value = 1588
$ws4.Range("F46").Value = 239
$ws4.Range("F47").Value = 84
$ws4.Range("F48").Value = 388
$ws4.Range("F50").Value = 3518
